# The "Exclude Table" header previously spanned two merged-looking cells
# (L1 = "Exclude Table", M1 = "Exclude"). The sheet is simplified to a
# single header cell: L1 becomes "#Exclude Table" and M1 is cleared out,
# which also removes the now-unused "Exclude Table"/"Exclude" strings
# from the shared string table (replaced by the single "#Exclude Table").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "#Exclude Table"
$ws.Range("M1").ClearContents()

# Move the active selection to L1 (was M9 before the edit).
$ws.Range("L1").Select()
